$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.113.99"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "2.739.84"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.97"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.37"
$ws.Range("E6").Value = "  +6.58%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").Value = "2.738.76"
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("E11").Value = "  +6.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.95"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").Value = "3.238.07"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("D17").Value = "68.934.40"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "2.722.08"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.97"
$ws.Range("E19").Value = "  +5.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "376.84"
$ws.Range("E20").Value = "  +5.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.72"
$ws.Range("E21").Value = "  +5.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.57"
$ws.Range("E22").Value = "  +4.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.02"
$ws.Range("E23").Value = "  +6.47%  "
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.05"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +4.57%  "
$ws.Range("D28").Value = "2.872.43"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "590.97"
$ws.Range("E30").Value = "  +6.40%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.39"
$ws.Range("E32").Value = "  +5.50%  "
$ws.Range("E33").Value = "  +6.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.98"
$ws.Range("E34").Value = "  +6.30%  "
$ws.Range("E35").Value = "  +5.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.63"
$ws.Range("E36").Value = "  +3.85%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.13"
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.42"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.384"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("E41").Value = "  +4.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.53"
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").Value = "  +4.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.01"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "0.0₆0311"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.11"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.80"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.99"
$ws.Range("E49").Value = "  +5.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.81"
$ws.Range("E50").Value = "  +8.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.608"
$ws.Range("E51").Value = "  +7.30%  "
